# feat: add 2022-Q4 data
#
# - Insert a new "2022-Q4" sheet right after "总计" and before "2022-Q3"
#   (so the sheet order becomes 总计, 2022-Q4, 2022-Q3, 2022-Q2), populated
#   with the Q4 fund-position data.
# - Update the "总计" (summary) sheet so it lists 2022-Q4 first, followed
#   by the pre-existing 2022-Q3 / 2022-Q2 rows (shifted down by one row).
#
# NOTE: worksheet references obtained via Worksheets.Item(...) are
# positional snapshots — once Worksheets.Add() inserts a new sheet the
# sheet collection shifts and old variables can silently start pointing
# at a different sheet. To stay safe, sheet references are re-fetched by
# name immediately before each use, after the Add() call.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet, positioned right after "总计".
# ---------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q4Sheet.Name = "2022-Q4"

# Clone the header row + styling (including the bordered/bold header
# style and the index-column style) from the "2022-Q3" sheet so the new
# sheet matches the workbook's existing look, then extend the
# 2-data-row template down to the 4 rows "2022-Q4" needs.
# (Re-fetch "2022-Q3" / "2022-Q4" fresh here — the sheet collection just
# shifted because of Add() above.)
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q4Sheet = $wb.Worksheets.Item("2022-Q4")

$q3Sheet.Range("B1:H3").Copy($q4Sheet.Range("B1"))
$q3Sheet.Range("A2:A3").Copy($q4Sheet.Range("A2"))
$q3Sheet.Range("A3:H3").Copy($q4Sheet.Range("A4"))
$q3Sheet.Range("A3:H3").Copy($q4Sheet.Range("A5"))

# Fund codes / scale / position / ratio / market-value columns are
# text-like (e.g. leading-zero fund codes "014277"), so force text
# NumberFormat BEFORE assigning their values — otherwise numeric-looking
# content gets silently coerced into actual numbers (losing the leading
# zero, trailing zeros, etc). Column A (index) / H (rank) stay numeric.
$q4Sheet.Range("B2:B5").NumberFormat = "@"
$q4Sheet.Range("D2:G5").NumberFormat = "@"

# Fill in the actual 2022-Q4 values.
$q4Data = @(
    @(0, "014277", "万家北交所慧选两年定期开放混合A", "3.25", "94.43", "8.48", "0.2756", 1),
    @(1, "016307", "景顺长城北交所精选两年定开混合A", "1.83", "43.56", "2.77", "0.0507", 4),
    @(2, "014278", "万家北交所慧选两年定期开放混合C", "0.45", "94.43", "8.48", "0.0382", 1),
    @(3, "016308", "景顺长城北交所精选两年定开混合C", "0.27", "43.56", "2.77", "0.0075", 4)
)

$r = 2
foreach ($row in $q4Data) {
    $q4Sheet.Cells.Item($r, 1).Value = $row[0]
    $q4Sheet.Cells.Item($r, 2).Value = $row[1]
    $q4Sheet.Cells.Item($r, 3).Value = $row[2]
    $q4Sheet.Cells.Item($r, 4).Value = $row[3]
    $q4Sheet.Cells.Item($r, 5).Value = $row[4]
    $q4Sheet.Cells.Item($r, 6).Value = $row[5]
    $q4Sheet.Cells.Item($r, 7).Value = $row[6]
    $q4Sheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: 2022-Q4 becomes the first data row, with
#    the old 2022-Q3 / 2022-Q2 rows pushed down by one.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Use the existing A2 index-cell formatting as the template for the new
# A4 cell (same bold/bordered style as A2 / A3).
$totalSheet.Cells.Item(2, 1).Copy($totalSheet.Cells.Item(4, 1))

$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(4, 2).Value = "2022-Q2"
$totalSheet.Cells.Item(4, 3).Value = 2
$totalSheet.Cells.Item(4, 4).Value = 0.18

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(3, 3).Value = 2
$totalSheet.Cells.Item(3, 4).Value = 0.21

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 4
$totalSheet.Cells.Item(2, 4).Value = 0.37

# ---------------------------------------------------------------------
# 3. Keep "2022-Q2" as the active/selected tab, matching the original
#    workbook's selection state.
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Activate()
